$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.669.97"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = "'1.800.45"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('D5').Value = "'230.91"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').Value = "'0.5940"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = "'0.2766"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').Value = "'0.06803"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.42%  '
$ws.Range('D10').Value = "'23.31"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('D11').Value = "'0.07518"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').Value = "'1.802.61"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').Value = "'4.687"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.76%  '
$ws.Range('D14').Value = "'0.6250"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.98%  '
$ws.Range('D15').Value = "'2.045.05"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = "'0.000009123"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.03%  '
$ws.Range('D17').Value = "'75.30"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.22%  '
$ws.Range('D18').Value = "'28.561.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').Value = "'5.445"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.93%  '
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').Value = "'209.92"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.29%  '
$ws.Range('D22').Value = "'11.39"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('E23').Value = '  -3.10%  '
$ws.Range('D24').Value = "'1.004"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').Value = "'154.27"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').Value = "'7.829"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.40%  '
$ws.Range('D27').Value = "'0.1272"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = "'16.39"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('D29').Value = "'1.447"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('D30').Value = "'0.06369"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('D31').Value = "'1.415"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.79%  '
$ws.Range('D32').Value = "'3.739"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('D33').Value = "'3.713"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').Value = "'1.703"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.23%  '
$ws.Range('D35').Value = "'1.048"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.32%  '
$ws.Range('D36').Value = "'0.6337"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.21%  '
$ws.Range('D37').Value = "'2.507"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('D38').Value = "'2.712"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').Value = "'0.01698"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').Value = "'6.378"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.09%  '
$ws.Range('D41').Value = "'1.134.32"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.44%  '
$ws.Range('D42').Value = "'0.8627"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.74%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').Value = "'100.66"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = "'1.963.91"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = "'60.38"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.03%  '
$ws.Range('E47').Value = '  -5.46%  '
$ws.Range('D48').Value = "'1.577"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').Value = "'8.330"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').Value = "'0.4500"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').Value = "'0.05441"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.75%  '
